# Add a new weekly price record (row 96) to the sheet, mirroring the
# structure of the existing rows (e.g. row 95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96

$ws.Range("A$row").Value2 = 11
$ws.Range("B$row").Value2 = "Vega Monumental Concepción"
$ws.Range("C$row").Value2 = "Bíobío"
$ws.Range("D$row").Value2 = 44656
$ws.Range("D$row").NumberFormat = $ws.Range("D95").NumberFormat
$ws.Range("E$row").Value2 = 8
$ws.Range("F$row").Value2 = "Fruta"
$ws.Range("G$row").Value2 = 100101
$ws.Range("H$row").Value2 = "Berries"
$ws.Range("I$row").Value2 = 100101001
$ws.Range("J$row").Value2 = "Arándano (blue)"
$ws.Range("K$row").Value2 = "Sin especificar"
$ws.Range("L$row").Value2 = "Primera"
$ws.Range("M$row").Value2 = 200
$ws.Range("N$row").Value2 = 3500
$ws.Range("O$row").Value2 = 4000
$ws.Range("P$row").Value2 = 3750
$ws.Range("Q$row").Value2 = "$/bandeja 2 kilos"
$ws.Range("R$row").Value2 = "Región de Ñuble"
$ws.Range("S$row").Value2 = 1875
$ws.Range("T$row").Value2 = 2
